# Updated symbol list on Tue Dec 27 10:28:11 UTC 2022 with GitHub Actions
# This script refreshes the "Price" column (and two "Worstin24h" tags that
# moved in the source feed) on Sheet1 to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (cell, new value). Price cells are numeric-looking text in
# the source data (e.g. "0.05960" must keep its trailing zero), so we force
# the cell to Text format before writing the value - otherwise Excel would
# coerce the literal into a real number and lose formatting/precision.
$priceUpdates = @(
    @("D2",  "243.83"),
    @("D3",  "23.07"),
    @("D4",  "5.386"),
    @("D5",  "0.05960"),
    @("D6",  "3.430"),
    @("D7",  "6.512"),
    @("D8",  "0.8110"),
    @("D9",  "0.9262"),
    @("D10", "0.1431"),
    @("D11", "0.07411"),
    @("D12", "0.03277"),
    @("D13", "0.03090"),
    @("D14", "0.09355"),
    @("D15", "3.863"),
    @("D16", "0.001587"),
    @("D17", "0.04697"),
    @("D18", "0.0005949"),
    @("D19", "0.005940"),
    @("D20", "0.001262"),
    @("D21", "0.004795"),
    @("D22", "0.00007997"),
    @("D23", "3.573"),
    @("D26", "0.1331"),
    @("D27", "0.0002339"),
    @("D40", "0.03933"),
    @("D41", "0.006378"),
    @("D42", "0.1078"),
    @("D43", "0.002579"),
    @("D44", "0.008968"),
    @("D45", "0.00005182"),
    @("D47", "0.6849"),
    @("D48", "0.002145"),
    @("D49", "0.00002100"),
    @("D50", "0.0002000")
)

# An untouched, never-referenced cell whose style is the workbook default -
# used to snap each edited cell's style back to "no explicit style" after
# the temporary Text formatting has done its job of protecting the literal.
$defaultStyleCell = $ws.Range("Z100")

foreach ($pair in $priceUpdates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $defaultStyleCell.Style
}

# The "Worstin24h" marker moved from row 18 (One/ONE) to row 27 (UpBots/UBXT).
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E27").Value = "26UpBotsUBXT"
